# Updates the cryptos list data (price / volume / coin-name swaps) as part
# of the scheduled "Updated cryptos list ... with GitHub Actions" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, ColumnLetter, NewValue
$updates = @(
    @(2,  "D", "58.037.46"),
    @(2,  "E", "  -0.51%  "),

    @(3,  "D", "2.348.18"),
    @(3,  "E", "  +0.68%  "),

    @(4,  "E", "  -0.42%  "),

    @(5,  "D", "541.68"),
    @(5,  "E", "  -0.44%  "),

    @(6,  "D", "134.17"),
    @(6,  "E", "  -0.69%  "),

    @(7,  "E", "  -0.01%  "),

    @(8,  "E", "  +5.31%  "),

    @(9,  "E", "  +1.42%  "),

    @(10, "D", "5.52"),
    @(10, "E", "  +2.27%  "),

    @(11, "E", "  -1.95%  "),

    @(12, "D", "0.356"),
    @(12, "E", "  -0.27%  "),

    # Rows 13 and 14 swap: Avalanche moves to row 13, WrappedliquidstakedEther2.0 to row 14
    @(13, "B", "Avalanche"),
    @(13, "C", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"),
    @(13, "D", "23.77"),
    @(13, "E", "  +0.45%  "),

    @(14, "B", "WrappedliquidstakedEther2.0"),
    @(14, "C", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"),
    @(14, "D", "2.764.51"),
    @(14, "E", "  -0.73%  "),

    @(15, "D", "57.962.86"),
    @(15, "E", "  -0.54%  "),

    @(16, "E", "  +0.70%  "),

    @(17, "D", "2.337.69"),
    @(17, "E", "  -1.39%  "),

    @(18, "D", "10.74"),

    @(19, "E", "  +1.98%  "),

    @(20, "D", "328.71"),
    @(20, "E", "  -2.12%  "),

    @(21, "D", "6.73"),
    @(21, "E", "  -0.20%  "),

    @(22, "E", "  +0.14%  "),

    @(23, "E", "  +1.76%  "),

    @(24, "D", "0.164"),
    @(24, "E", "  -3.36%  "),

    @(25, "D", "0.995"),
    @(25, "E", "  -0.42%  "),

    @(26, "D", "8.38"),
    @(26, "E", "  -1.15%  "),

    @(27, "D", "1.32"),
    @(27, "E", "  -7.69%  "),

    @(28, "E", "  -0.11%  "),

    @(29, "D", "170.35"),
    @(29, "E", "  -0.12%  "),

    @(30, "E", "  -0.29%  "),

    @(31, "E", "  -0.27%  "),

    @(32, "D", "18.33"),
    @(32, "E", "  -1.29%  "),

    @(33, "E", "  -0.98%  "),

    @(34, "E", "  -0.02%  "),

    @(35, "E", "  +0.02%  "),

    @(36, "D", "4.19"),
    @(36, "E", "  -0.30%  "),

    @(37, "E", "  -3.02%  "),

    @(38, "D", "1.60"),
    @(38, "E", "  -2.38%  "),

    @(39, "D", "39.13"),
    @(39, "E", "  -0.59%  "),

    @(40, "D", "0.379"),
    @(40, "E", "  -0.49%  "),

    @(41, "D", "290.03"),
    @(41, "E", "  +1.02%  "),

    @(42, "D", "140.53"),
    @(42, "E", "  -6.68%  "),

    @(43, "D", "3.64"),
    @(43, "E", "  +0.00%  "),

    @(44, "E", "  +2.19%  "),

    @(45, "E", "  +0.60%  "),

    # Rows 46 and 47 swap: InjectiveProtocol moves to row 46, Mantle to row 47
    @(46, "B", "InjectiveProtocol"),
    @(46, "C", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"),
    @(46, "D", "18.91"),
    @(46, "E", "  -2.20%  "),

    @(47, "B", "Mantle"),
    @(47, "C", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"),
    @(47, "D", "0.566"),
    @(47, "E", "  +0.21%  "),

    @(48, "E", "  +1.32%  "),

    @(49, "E", "  +0.17%  "),

    @(50, "D", "11.07"),
    @(50, "E", "  +0.07%  "),

    @(51, "D", "4.70"),
    @(51, "E", "  +0.66%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Range("$col$row").Value = $val
}
